$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change existing row 2, column B value from "2C#2" to "7C#14"
$ws.Range("B2").Value = "7C#14"

# Add a new row 3 duplicating the original row 2 data (Pull # 2)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "2C#2"
$ws.Range("C3").Value = "EXPRESS"
$ws.Range("D3").Value = "100+00"
$ws.Range("E3").Value = "200+00"

# Update the selection to match the new active cell
$ws.Range("E7").Select()
